$d = $word.ActiveDocument

$replacements = @(
    @("49×36=", "37×56="),
    @("36×90=", "13×76="),
    @("16×33=", "96×38="),
    @("37×12=", "11×31="),
    @("36×81=", "61×81="),
    @("95×70=", "97×57="),
    @("11×85=", "25×29="),
    @("57×80=", "16×12="),
    @("85×27=", "31×64="),
    @("93×54=", "34×89="),
    @("25×98=", "89×16="),
    @("46×34=", "12×15="),
    @("45×13=", "86×88="),
    @("57×32=", "81×95="),
    @("77×27=", "67×73="),
    @("40×60=", "61×47="),
    @("22×98=", "59×14="),
    @("96×40=", "59×84="),
    @("99×35=", "55×49="),
    @("16×85=", "55×83="),
    @("60×49=", "68×19="),
    @("39×49=", "79×36="),
    @("69×21=", "31×28="),
    @("70×89=", "24×73="),
    @("56×80=", "66×49=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
